# Trade #85 (momentum strategy, "Downward momentum" setup opened 00:21:03)
# closed early at 2026-02-18 00:29:01 with ~+0.000% (essentially flat) P&L,
# and a brand-new MarketMaking trade (#143) was opened at 00:28:55.
# This script reproduces every cell touched by that commit.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet - roll-up numbers after the close
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1499.4    # Current Capital
$summary.Range("B4").Value = 0.51      # Total P&L $
$summary.Range("B6").Value = 113       # Total Trades
$summary.Range("B7").Value = 54        # Winning Trades
$summary.Range("B9").Value = 47.79     # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status sheet - momentum strategy row (row 11)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C11").Value = 99.23
$status.Range("D11").Value = 31
$status.Range("E11").Value = -0.76
$status.Range("F11").Value = -0.77
$status.Range("G11").Value = 29.03

# ---------------------------------------------------------------------------
# All Trades sheet - trade #113 (row 114) closes out
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Range("G114").Value = 0.97
$allTrades.Range("H114").Value = "CLOSED"
$allTrades.Range("I114").Value = 1.0417
$allTrades.Range("J114").Value = 0.01
$allTrades.Range("K114").Value = 99.23
$allTrades.Range("L114").Value = "early_exit"
$allTrades.Range("M114").Value = 0.17

# All Trades sheet - new trade #143 (row 144), MarketMaking, still OPEN
$allTrades.Range("A144").Value = 143
$allTrades.Range("B144").Value = "'2026-02-18"
$allTrades.Range("B144").Style = "Normal"
$allTrades.Range("C144").Value = "00:28:55"
$allTrades.Range("D144").Value = "MarketMaking"
$allTrades.Range("E144").Value = "DOWN"
$allTrades.Range("F144").Value = 0.96
$allTrades.Range("G144").Value = "'"
$allTrades.Range("G144").Style = "Normal"
$allTrades.Range("H144").Value = "OPEN"
$allTrades.Range("I144").Value = 0
$allTrades.Range("J144").Value = 0
$allTrades.Range("K144").Value = 99.47967800952271
$allTrades.Range("L144").Value = "'"
$allTrades.Range("L144").Style = "Normal"
$allTrades.Range("M144").Value = 0
$allTrades.Range("N144").Value = 0
$allTrades.Range("O144").Value = 0
$allTrades.Range("P144").Value = 0.6
$allTrades.Range("Q144").Value = "Normal spread capture: 190 bps"

# ---------------------------------------------------------------------------
# momentum sheet - its own copy of trade #113 (row 32) closes out too
# ---------------------------------------------------------------------------
$momentum = $wb.Worksheets.Item("momentum")
$momentum.Range("G32").Value = 0.97
$momentum.Range("H32").Value = "CLOSED"
$momentum.Range("I32").Value = 1.0417
$momentum.Range("J32").Value = 0.01
$momentum.Range("K32").Value = 99.23
$momentum.Range("P32").Value = "early_exit"
$momentum.Range("Q32").Value = 0.17

# ---------------------------------------------------------------------------
# MarketMaking sheet - its own copy of the new trade #143 (row 59)
# ---------------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
$marketMaking.Range("A59").Value = 143
$marketMaking.Range("B59").Value = "'2026-02-18"
$marketMaking.Range("B59").Style = "Normal"
$marketMaking.Range("C59").Value = "00:28:55"
$marketMaking.Range("D59").Value = "MarketMaking"
$marketMaking.Range("E59").Value = "DOWN"
$marketMaking.Range("F59").Value = 0.96
$marketMaking.Range("G59").Value = "'"
$marketMaking.Range("G59").Style = "Normal"
$marketMaking.Range("H59").Value = "OPEN"
$marketMaking.Range("I59").Value = 0
$marketMaking.Range("J59").Value = 0
$marketMaking.Range("K59").Value = 99.47967800952271
$marketMaking.Range("L59").Value = 0
$marketMaking.Range("M59").Value = 0
$marketMaking.Range("N59").Value = 0.6
$marketMaking.Range("O59").Value = "Normal spread capture: 190 bps"
$marketMaking.Range("P59").Value = "'"
$marketMaking.Range("P59").Style = "Normal"
$marketMaking.Range("Q59").Value = 0

Write-Host "Applied trade #85 close-out and trade #143 open edits"
